# Auto-generated Excel COM-interop script applying scheduled runner updates
# to the per-leve market-price columns (H-N) across all sheets.
$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

$ALC.Range("H12").Value = 127.77273
$ALC.Range("I12").Value = 126.71429
$ALC.Range("K12").Value = 126.71429
$ALC.Range("M12").Value = 43.28570999999999
$ALC.Range("H43").Value = 2495.5
$ALC.Range("I43").Value = 2494.3333
$ALC.Range("J43").Value = 2499
$ALC.Range("K43").Value = 2494.3333
$ALC.Range("L43").Value = 2499
$ALC.Range("M43").Value = -2425.3333
$ALC.Range("N43").Value = -2637
$ALC.Range("H74").Value = 7777.773
$ALC.Range("I74").Value = 3500.6
$ALC.Range("J74").Value = 9035.764999999999
$ALC.Range("K74").Value = 3500.6
$ALC.Range("L74").Value = 9035.764999999999
$ALC.Range("M74").Value = -2564.6
$ALC.Range("N74").Value = -10907.765
$ALC.Range("H77").Value = 7777.773
$ALC.Range("I77").Value = 3500.6
$ALC.Range("J77").Value = 9035.764999999999
$ALC.Range("K77").Value = 17503
$ALC.Range("L77").Value = 45178.825
$ALC.Range("M77").Value = -12823
$ALC.Range("N77").Value = -54538.825
$ALC.Range("H86").Value = 1585.3
$ALC.Range("I86").Value = 1044.909
$ALC.Range("J86").Value = 2245.7778
$ALC.Range("K86").Value = 1044.909
$ALC.Range("L86").Value = 2245.7778
$ALC.Range("M86").Value = 78.09099999999989
$ALC.Range("N86").Value = -4491.7778
$ALC.Range("H89").Value = 1585.3
$ALC.Range("I89").Value = 1044.909
$ALC.Range("J89").Value = 2245.7778
$ALC.Range("K89").Value = 5224.545
$ALC.Range("L89").Value = 11228.889
$ALC.Range("M89").Value = 391.4549999999999
$ALC.Range("N89").Value = -22460.889
$ALC.Range("H92").Value = 217.44444
$ALC.Range("I92").Value = 195.35715
$ALC.Range("J92").Value = 294.75
$ALC.Range("K92").Value = 195.35715
$ALC.Range("L92").Value = 294.75
$ALC.Range("M92").Value = 1052.64285
$ALC.Range("N92").Value = -2790.75
$ALC.Range("H94").Value = 11342.714
$ALC.Range("I94").Value = 4899.8335
$ALC.Range("K94").Value = 4899.8335
$ALC.Range("M94").Value = -4448.8335
$ALC.Range("H98").Value = 1532.9333
$ALC.Range("I98").Value = 1550.6296
$ALC.Range("K98").Value = 1550.6296
$ALC.Range("M98").Value = -52.62959999999998
$ALC.Range("H112").Value = 5412.7847
$ALC.Range("J112").Value = 5469.2344
$ALC.Range("L112").Value = 16407.7032
$ALC.Range("N112").Value = -18623.7032
$ALC.Range("H122").Value = 1532.9333
$ALC.Range("I122").Value = 1550.6296
$ALC.Range("K122").Value = 4651.8888
$ALC.Range("M122").Value = -2201.8888
$ALC.Range("H138").Value = 2460.516
$ALC.Range("I138").Value = 1695.4
$ALC.Range("J138").Value = 2824.8572
$ALC.Range("K138").Value = 5086.200000000001
$ALC.Range("L138").Value = 8474.571599999999
$ALC.Range("M138").Value = 53.79999999999927
$ALC.Range("N138").Value = -18754.5716
$ARM.Range("H32").Value = 20917.072
$ARM.Range("I32").Value = 21297.203
$ARM.Range("K32").Value = 21297.203
$ARM.Range("M32").Value = -21010.203
$ARM.Range("H61").Value = 2223015.8
$ARM.Range("I61").Value = 2381758.8
$ARM.Range("J61").Value = 614
$ARM.Range("K61").Value = 2381758.8
$ARM.Range("L61").Value = 614
$ARM.Range("M61").Value = -2381546.8
$ARM.Range("N61").Value = -1038
$ARM.Range("H74").Value = 3339.8518
$ARM.Range("I74").Value = 1013.05884
$ARM.Range("J74").Value = 7295.4
$ARM.Range("K74").Value = 1013.05884
$ARM.Range("L74").Value = 7295.4
$ARM.Range("M74").Value = -139.05884
$ARM.Range("N74").Value = -9043.4
$ARM.Range("H77").Value = 3339.8518
$ARM.Range("I77").Value = 1013.05884
$ARM.Range("J77").Value = 7295.4
$ARM.Range("K77").Value = 5065.2942
$ARM.Range("L77").Value = 36477
$ARM.Range("M77").Value = -697.2942000000003
$ARM.Range("N77").Value = -45213
$ARM.Range("H97").Value = 1079.6428
$ARM.Range("I97").Value = 893.4
$ARM.Range("K97").Value = 893.4
$ARM.Range("M97").Value = -397.4
$ARM.Range("H102").Value = 28715.23
$ARM.Range("J102").Value = 594.5
$ARM.Range("L102").Value = 594.5
$ARM.Range("N102").Value = -3838.5
$ARM.Range("H122").Value = 2140.4
$ARM.Range("I122").Value = 2140.4
$ARM.Range("K122").Value = 6421.200000000001
$ARM.Range("M122").Value = -3971.200000000001
$ARM.Range("H136").Value = 2223015.8
$ARM.Range("I136").Value = 2381758.8
$ARM.Range("J136").Value = 614
$ARM.Range("K136").Value = 7145276.399999999
$ARM.Range("L136").Value = 1842
$ARM.Range("M136").Value = -7142726.399999999
$ARM.Range("N136").Value = -6942
$BSM.Range("H86").Value = 1257.1613
$BSM.Range("I86").Value = 1257.1613
$BSM.Range("J86").Value = 0
$BSM.Range("K86").Value = 1257.1613
$BSM.Range("L86").Value = 0
$BSM.Range("M86").Value = -134.1613
$BSM.Range("N86").ClearContents()
$BSM.Range("H89").Value = 1257.1613
$BSM.Range("I89").Value = 1257.1613
$BSM.Range("J89").Value = 0
$BSM.Range("K89").Value = 6285.8065
$BSM.Range("L89").Value = 0
$BSM.Range("M89").Value = -669.8064999999997
$BSM.Range("N89").ClearContents()
$BSM.Range("H94").Value = 1452
$BSM.Range("I94").Value = 1220.25
$BSM.Range("J94").Value = 2070
$BSM.Range("K94").Value = 1220.25
$BSM.Range("L94").Value = 2070
$BSM.Range("M94").Value = -769.25
$BSM.Range("N94").Value = -2972
$BSM.Range("H105").Value = 4087.087
$BSM.Range("I105").Value = 3975.8572
$BSM.Range("K105").Value = 3975.8572
$BSM.Range("M105").Value = -2228.8572
$BSM.Range("H131").Value = 68998.5
$BSM.Range("J131").Value = 68998.5
$BSM.Range("L131").Value = 68998.5
$BSM.Range("N131").Value = -79078.5
$CRP.Range("H16").Value = 816.5
$CRP.Range("I16").Value = 588.6667
$CRP.Range("K16").Value = 588.6667
$CRP.Range("M16").Value = -301.6667
$CRP.Range("H31").Value = 11270.733
$CRP.Range("J31").Value = 27294.428
$CRP.Range("L31").Value = 27294.428
$CRP.Range("N31").Value = -27884.428
$CRP.Range("H34").Value = 11270.733
$CRP.Range("J34").Value = 27294.428
$CRP.Range("L34").Value = 27294.428
$CRP.Range("N34").Value = -27698.428
$CRP.Range("H58").Value = 888059.1
$CRP.Range("I58").Value = 1374202.1
$CRP.Range("K58").Value = 1374202.1
$CRP.Range("M58").Value = -1373999.1
$CRP.Range("H93").Value = 34165.332
$CRP.Range("I93").Value = 34165.332
$CRP.Range("K93").Value = 34165.332
$CRP.Range("M93").Value = -32293.332
$CRP.Range("H99").Value = 3112.4783
$CRP.Range("I99").Value = 2841.2727
$CRP.Range("J99").Value = 3361.0833
$CRP.Range("K99").Value = 2841.2727
$CRP.Range("L99").Value = 3361.0833
$CRP.Range("M99").Value = -1343.2727
$CRP.Range("N99").Value = -6357.0833
$CRP.Range("H105").Value = 33245.453
$CRP.Range("I105").Value = 44587.75
$CRP.Range("K105").Value = 44587.75
$CRP.Range("M105").Value = -42840.75
$CRP.Range("H113").Value = 816.5
$CRP.Range("I113").Value = 588.6667
$CRP.Range("K113").Value = 588.6667
$CRP.Range("M113").Value = 1581.3333
$CRP.Range("H126").Value = 3112.4783
$CRP.Range("I126").Value = 2841.2727
$CRP.Range("J126").Value = 3361.0833
$CRP.Range("K126").Value = 8523.8181
$CRP.Range("L126").Value = 10083.2499
$CRP.Range("M126").Value = -6053.8181
$CRP.Range("N126").Value = -15023.2499
$CRP.Range("H136").Value = 888059.1
$CRP.Range("I136").Value = 1374202.1
$CRP.Range("K136").Value = 4122606.3
$CRP.Range("M136").Value = -4120056.3
$CUL.Range("H32").Value = 390476800
$CUL.Range("I32").Value = 1249.5
$CUL.Range("K32").Value = 3748.5
$CUL.Range("M32").Value = -3465.5
$GSM.Range("H2").Value = 938.2
$GSM.Range("I2").Value = 1440.875
$GSM.Range("J2").Value = 363.7143
$GSM.Range("K2").Value = 1440.875
$GSM.Range("L2").Value = 363.7143
$GSM.Range("M2").Value = -1327.875
$GSM.Range("N2").Value = -589.7143
$GSM.Range("H97").Value = 894.3461
$GSM.Range("I97").Value = 917.5909
$GSM.Range("K97").Value = 917.5909
$GSM.Range("M97").Value = -421.5909
$GSM.Range("H139").Value = 145000
$GSM.Range("J139").Value = 145000
$GSM.Range("L139").Value = 145000
$GSM.Range("N139").Value = -155280
$LTW.Range("H93").Value = 2573.75
$LTW.Range("I93").Value = 2573.75
$LTW.Range("K93").Value = 2573.75
$LTW.Range("M93").Value = -1325.75
$LTW.Range("H132").Value = 3487368.5
$LTW.Range("J132").Value = 9749.75
$LTW.Range("L132").Value = 29249.25
$LTW.Range("N132").Value = -34309.25
$WVR.Range("H81").Value = 4379.0835
$WVR.Range("I81").Value = 4256.25
$WVR.Range("K81").Value = 8512.5
$WVR.Range("M81").Value = -7451.5
$WVR.Range("H84").Value = 4379.0835
$WVR.Range("I84").Value = 4256.25
$WVR.Range("K84").Value = 42562.5
$WVR.Range("M84").Value = -37258.5
$WVR.Range("H107").Value = 1382.6923
$WVR.Range("I107").Value = 813
$WVR.Range("J107").Value = 2294.2
$WVR.Range("K107").Value = 2439
$WVR.Range("L107").Value = 6882.599999999999
$WVR.Range("M107").Value = -519
$WVR.Range("N107").Value = -10722.6
